$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observation record appended as row 26
$ws.Range("A26").Value = 104165878
$ws.Range("B26").Value = 5207
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 100155
$ws.Range("F26").Value = "Mindre timmerman"
$ws.Range("G26").Value = "Acanthocinus griseus"
$ws.Range("H26").Value = "(Fabricius, 1792)"

$ws.Range("P26").Value = "Finsjöbrännan, Sm"
$ws.Range("Q26").Value = 575416
$ws.Range("R26").Value = 6336378
$ws.Range("S26").Value = 25
$ws.Range("T26").Value = "Kalmar"
$ws.Range("U26").Value = "Mönsterås"
$ws.Range("V26").Value = "Småland"
$ws.Range("W26").Value = "Fliseryd"

# Startdatum / Slutdatum are stored as plain text ("2022-06-10"), not as a
# date serial. Enter them with a leading apostrophe so they are kept as
# literal text, then clear the resulting cell format so no extra
# number-format/quote-prefix styling is left behind on the cell.
$ws.Range("Y26").Value = "'2022-06-10"
$ws.Range("Y26").ClearFormats()
$ws.Range("AA26").Value = "'2022-06-10"
$ws.Range("AA26").ClearFormats()

$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false

$ws.Range("AW26").Value = "Olof Persson"
$ws.Range("AX26").Value = "Olof Persson, Jesper Hansson"

# The remaining columns of the new row are present in the source data but
# hold no value. Touch them with a no-op formatting change so a (blank)
# cell is recorded at each of these positions, matching the shape of the
# appended row.
$blankCols = @("I26","J26","K26","L26","M26","N26","AF26","AT26","AY26")
foreach ($addr in $blankCols) {
    $ws.Range($addr).Font.Bold = $false
}
